$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $style = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $style
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "66.510.83"
$ws.Range("E2").Value = "  +4.26%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.485.07"
$ws.Range("E3").Value = "  +2.70%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "588.71"
$ws.Range("E5").Value = "  +3.17%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "168.74"
$ws.Range("E6").Value = "  +4.38%  "

# Row 8 - LidoStakedEther
Set-TextValue $ws.Range("D8") "3.481.29"
$ws.Range("E8").Value = "  +2.54%  "

# Row 9 - XRP
Set-TextValue $ws.Range("D9") "0.591"
$ws.Range("E9").Value = "  +7.78%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +6.54%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("D12") "0.437"
$ws.Range("E12").Value = "  +3.85%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "4.087.36"
$ws.Range("E13").Value = "  +2.72%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.00%  "

# Row 15 - Avalanche
Set-TextValue $ws.Range("D15") "28.13"
$ws.Range("E15").Value = "  +4.63%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +3.66%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "66.538.72"
$ws.Range("E17").Value = "  +4.21%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "3.488.69"
$ws.Range("E18").Value = "  +3.02%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "6.32"
$ws.Range("E19").Value = "  +3.54%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +3.87%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "391.81"
$ws.Range("E21").Value = "  +4.30%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +2.01%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "73.11"
$ws.Range("E23").Value = "  +4.39%  "

# Row 24 - Dai
Set-TextValue $ws.Range("D24") "0.999"
$ws.Range("E24").Value = "  -0.18%  "

# Row 25 - Polygon
Set-TextValue $ws.Range("D25") "0.535"
$ws.Range("E25").Value = "  +4.31%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +5.63%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  +7.51%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +1.99%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  +0.04%  "

# Row 30 - NEARProtocol
$ws.Range("E30").Value = "  +4.26%  "

# Row 31 - Fetch.AI
Set-TextValue $ws.Range("D31") "1.46"
$ws.Range("E31").Value = "  +5.94%  "

# Row 32 - PancakeSwap
Set-TextValue $ws.Range("D32") "2.06"
$ws.Range("E32").Value = "  +3.01%  "

# Row 33 - EthereumClassic
Set-TextValue $ws.Range("D33") "23.58"
$ws.Range("E33").Value = "  +3.45%  "

# Row 34 - Aptos
$ws.Range("E34").Value = "  +4.88%  "

# Row 35 - USDe
$ws.Range("E35").Value = "  +0.01%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +9.16%  "

# Row 37 - Monero
Set-TextValue $ws.Range("D37") "161.79"
$ws.Range("E37").Value = "  +1.22%  "

# Row 38 - Mantle
$ws.Range("E38").Value = "  +3.28%  "

# Row 39 - Stacks
Set-TextValue $ws.Range("D39") "1.92"
$ws.Range("E39").Value = "  +6.46%  "

# Row 40 - RenderToken
$ws.Range("E40").Value = "  +4.54%  "

# Row 41 - was Filecoin, now Hedera
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D41") "0.0743"
$ws.Range("E41").Value = "  +3.19%  "

# Row 42 - was Hedera, now Filecoin
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D42") "4.64"
$ws.Range("E42").Value = "  +6.33%  "

# Row 43 - EnergySwap
Set-TextValue $ws.Range("D43") "26.51"
$ws.Range("E43").Value = "  +3.53%  "

# Row 44 - InjectiveProtocol
Set-TextValue $ws.Range("D44") "27.00"
$ws.Range("E44").Value = "  +3.40%  "

# Row 45 - OKB
Set-TextValue $ws.Range("D45") "43.10"
$ws.Range("E45").Value = "  +0.74%  "

# Row 46 - Maker
Set-TextValue $ws.Range("D46") "2.777.11"
$ws.Range("E46").Value = "  +1.54%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  +1.54%  "

# Row 48 - dogwifhat
$ws.Range("E48").Value = "  +2.13%  "

# Row 49 - Bittensor
Set-TextValue $ws.Range("D49") "346.45"
$ws.Range("E49").Value = "  +6.13%  "

# Row 50 - ONDO
$ws.Range("E50").Value = "  +5.09%  "

# Row 51 - was Arweave, now SuiNetwork
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D51") "0.886"
$ws.Range("E51").Value = "  +9.39%  "
